$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price list refresh (prices in column D are stored as text,
# since some values use "." as a thousands separator, e.g. "51.626.53").
# Force text format before assigning so Excel does not reinterpret them
# as numbers (which would also lose formatting like "0.220").

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.626.53'
$ws.Range('E2').Value = '  -0.46%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.797.94'
$ws.Range('E3').Value = '  +0.56%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '354.61'
$ws.Range('E5').Value = '  -0.50%  '

$ws.Range('E6').Value = '  +0.78%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.559'
$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.623'
$ws.Range('E9').Value = '  +5.64%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.04'
$ws.Range('E10').Value = '  -0.30%  '

$ws.Range('E11').Value = '  +1.11%  '

$ws.Range('E12').Value = '  -0.86%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.77'
$ws.Range('E14').Value = '  +2.94%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.240.42'
$ws.Range('E15').Value = '  +0.35%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.800.15'
$ws.Range('E16').Value = '  +0.40%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.940'
$ws.Range('E17').Value = '  +0.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.599.73'
$ws.Range('E18').Value = '  -0.50%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.16'
$ws.Range('E20').Value = '  +2.48%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.32'
$ws.Range('E21').Value = '  +1.44%  '

$ws.Range('E22').Value = '  -0.31%  '

$ws.Range('E23').Value = '  +0.66%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.25'
$ws.Range('E24').Value = '  -0.51%  '

$ws.Range('E25').Value = '  +0.94%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.03'
$ws.Range('E27').Value = '  -1.70%  '

$ws.Range('E28').Value = '  +0.17%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.36'
$ws.Range('E29').Value = '  +0.77%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.51'

$ws.Range('E31').Value = '  +1.44%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.22'
$ws.Range('E32').Value = '  +8.69%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '52.24'
$ws.Range('E33').Value = '  +0.43%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.64'
$ws.Range('E34').Value = '  +8.82%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0446'
$ws.Range('E35').Value = '  -4.88%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0853'
$ws.Range('E36').Value = '  +1.20%  '

$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.74'
$ws.Range('E38').Value = '  -0.73%  '

$ws.Range('E39').Value = '  -1.95%  '

$ws.Range('E40').Value = '  +0.16%  '

$ws.Range('E41').Value = '  +0.35%  '

$ws.Range('E42').Value = '  -4.82%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.79'
$ws.Range('E43').Value = '  +0.16%  '

$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.19'
$ws.Range('E44').Value = '  -2.54%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.88'
$ws.Range('E45').Value = '  +0.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.140.35'
$ws.Range('E46').Value = '  +2.30%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.44'
$ws.Range('E47').Value = '  +5.64%  '

$ws.Range('E48').Value = '  +6.68%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.922'
$ws.Range('E49').Value = '  -3.21%  '

$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.38'
$ws.Range('E50').Value = '  +11.94%  '

$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.220'
$ws.Range('E51').Value = '  +16.15%  '
